{"js": "// Apply the arithmetic-expression text replacements described by the diff.\n// The document's single table holds one expression per cell (20 rows x 5 cols,\n// 100 cells total); the diff only changes each cell's <w:t> text, so we pair\n// each original expression with its replacement, in document (row-major) order,\n// and then write the whole table back via the `values` setter in one shot.\n// This sidesteps any risk of a later replacement's new text accidentally\n// matching an earlier replacement's new text during a text-search pass\n// (e.g. \"0+52=\" is a substring of \"40+52=\").\nconst replacements = [\n  [\"40+17=\", \"45+12=\"],\n  [\"55-48=\", \"87-47=\"],\n  [\"93-29=\", \"54+21=\"],\n  [\"9+83=\", \"45-43=\"],\n  [\"35+48=\", \"84-84=\"],\n  [\"31+24=\", \"92-89=\"],\n  [\"1+54=\", \"73-36=\"],\n  [\"45-34=\", \"49+5=\"],\n  [\"46+43=\", \"60+32=\"],\n  [\"81+16=\", \"59-11=\"],\n  [\"13+80=\", \"27+61=\"],\n  [\"65+17=\", \"58+39=\"],\n  [\"38+12=\", \"49-22=\"],\n  [\"51-19=\", \"50+44=\"],\n  [\"5+80=\", \"84+14=\"],\n  [\"57-27=\", \"46+37=\"],\n  [\"89-11=\", \"60-38=\"],\n  [\"40+51=\", \"34+10=\"],\n  [\"15+71=\", \"96-35=\"],\n  [\"38+3=\", \"6+11=\"],\n  [\"61-54=\", \"40+15=\"],\n  [\"51+10=\", \"12+18=\"],\n  [\"97-53=\", \"83-38=\"],\n  [\"53-2=\", \"66-46=\"],\n  [\"16-15=\", \"59+8=\"],\n  [\"67-60=\", \"88-46=\"],\n  [\"70+29=\", \"92-43=\"],\n  [\"4+7=\", \"79-3=\"],\n  [\"33+56=\", \"37+13=\"],\n  [\"90-88=\", \"76-40=\"],\n  [\"12+14=\", \"13+43=\"],\n  [\"68-26=\", \"45+46=\"],\n  [\"97-84=\", \"96-31=\"],\n  [\"20+65=\", \"16+45=\"],\n  [\"98-81=\", \"4+64=\"],\n  [\"3+47=\", \"91-47=\"],\n  [\"45+33=\", \"42+3=\"],\n  [\"85-71=\", \"8+64=\"],\n  [\"39+5=\", \"7+73=\"],\n  [\"49+1=\", \"86-58=\"],\n  [\"38+1=\", \"5+74=\"],\n  [\"59-57=\", \"27+69=\"],\n  [\"32+35=\", \"42+1=\"],\n  [\"93-80=\", \"46+1=\"],\n  [\"73-37=\", \"67-38=\"],\n  [\"62+15=\", \"46+35=\"],\n  [\"24+28=\", \"89-54=\"],\n  [\"96-88=\", \"31+39=\"],\n  [\"24+61=\", \"21+55=\"],\n  [\"84-47=\", \"32-25=\"],\n  [\"97-2=\", \"37+42=\"],\n  [\"24+1=\", \"0+35=\"],\n  [\"64-24=\", \"70-5=\"],\n  [\"46+47=\", \"17+60=\"],\n  [\"11+18=\", \"93-19=\"],\n  [\"41-35=\", \"39+45=\"],\n  [\"4+48=\", \"21+59=\"],\n  [\"59-54=\", \"41+11=\"],\n  [\"41-0=\", \"10-4=\"],\n  [\"97-68=\", \"63-32=\"],\n  [\"89-22=\", \"84-63=\"],\n  [\"22+34=\", \"44+20=\"],\n  [\"31+8=\", \"90-81=\"],\n  [\"62-27=\", \"67-53=\"],\n  [\"19+36=\", \"69-7=\"],\n  [\"99-10=\", \"81-74=\"],\n  [\"41+44=\", \"11+65=\"],\n  [\"78-63=\", \"58-14=\"],\n  [\"83-62=\", \"21+50=\"],\n  [\"39-32=\", \"41-36=\"],\n  [\"43+19=\", \"74-0=\"],\n  [\"3+21=\", \"65-10=\"],\n  [\"46-22=\", \"26-24=\"],\n  [\"29-16=\", \"41+57=\"],\n  [\"34-21=\", \"90-41=\"],\n  [\"86-62=\", \"86-29=\"],\n  [\"78+10=\", \"83-7=\"],\n  [\"63-31=\", \"20+38=\"],\n  [\"41+38=\", \"47-29=\"],\n  [\"94+1=\", \"15+50=\"],\n  [\"79+19=\", \"85-44=\"],\n  [\"14+15=\", \"49+11=\"],\n  [\"45-40=\", \"6+50=\"],\n  [\"86-14=\", \"72-27=\"],\n  [\"80-73=\", \"22+39=\"],\n  [\"44+35=\", \"62-32=\"],\n  [\"7+36=\", \"99-88=\"],\n  [\"91-61=\", \"74-46=\"],\n  [\"96-32=\", \"65-55=\"],\n  [\"59+16=\", \"24+35=\"],\n  [\"99-66=\", \"27+16=\"],\n  [\"63-51=\", \"40+52=\"],\n  [\"6+75=\", \"49+8=\"],\n  [\"19+61=\", \"61+26=\"],\n  [\"92-17=\", \"43-3=\"],\n  [\"24-20=\", \"52-16=\"],\n  [\"53-20=\", \"99-14=\"],\n  [\"39+40=\", \"53-14=\"],\n  [\"40+23=\", \"24+59=\"],\n  [\"0+52=\", \"21-3=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load('items');\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load('values');\nawait context.sync();\n\n// Build a lookup from the original expression text to its replacement,\n// then rewrite the table contents cell-by-cell, preserving layout/shape.\nconst lookup = new Map(replacements);\nconst newValues = table.values.map((row) =>\n  row.map((cellText) => {\n    const trimmed = cellText.trim();\n    return lookup.has(trimmed) ? lookup.get(trimmed) : cellText;\n  })\n);\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Apply the arithmetic-expression text replacements described by the diff.\n# The document's single table holds one expression per cell (20 rows x 5 cols,\n# 100 cells total); the diff only changes each cell's text, so we map each\n# original expression to its replacement and rewrite matching cells in place.\n# Using a content-keyed lookup (rather than blind positional indexing) means\n# the script is immune to any substring-collision between an old value and a\n# different cell's newly-written value (e.g. \"0+52=\" is a substring of \"40+52=\").\n$replacements = @{\n  '40+17=' = '45+12='\n  '55-48=' = '87-47='\n  '93-29=' = '54+21='\n  '9+83=' = '45-43='\n  '35+48=' = '84-84='\n  '31+24=' = '92-89='\n  '1+54=' = '73-36='\n  '45-34=' = '49+5='\n  '46+43=' = '60+32='\n  '81+16=' = '59-11='\n  '13+80=' = '27+61='\n  '65+17=' = '58+39='\n  '38+12=' = '49-22='\n  '51-19=' = '50+44='\n  '5+80=' = '84+14='\n  '57-27=' = '46+37='\n  '89-11=' = '60-38='\n  '40+51=' = '34+10='\n  '15+71=' = '96-35='\n  '38+3=' = '6+11='\n  '61-54=' = '40+15='\n  '51+10=' = '12+18='\n  '97-53=' = '83-38='\n  '53-2=' = '66-46='\n  '16-15=' = '59+8='\n  '67-60=' = '88-46='\n  '70+29=' = '92-43='\n  '4+7=' = '79-3='\n  '33+56=' = '37+13='\n  '90-88=' = '76-40='\n  '12+14=' = '13+43='\n  '68-26=' = '45+46='\n  '97-84=' = '96-31='\n  '20+65=' = '16+45='\n  '98-81=' = '4+64='\n  '3+47=' = '91-47='\n  '45+33=' = '42+3='\n  '85-71=' = '8+64='\n  '39+5=' = '7+73='\n  '49+1=' = '86-58='\n  '38+1=' = '5+74='\n  '59-57=' = '27+69='\n  '32+35=' = '42+1='\n  '93-80=' = '46+1='\n  '73-37=' = '67-38='\n  '62+15=' = '46+35='\n  '24+28=' = '89-54='\n  '96-88=' = '31+39='\n  '24+61=' = '21+55='\n  '84-47=' = '32-25='\n  '97-2=' = '37+42='\n  '24+1=' = '0+35='\n  '64-24=' = '70-5='\n  '46+47=' = '17+60='\n  '11+18=' = '93-19='\n  '41-35=' = '39+45='\n  '4+48=' = '21+59='\n  '59-54=' = '41+11='\n  '41-0=' = '10-4='\n  '97-68=' = '63-32='\n  '89-22=' = '84-63='\n  '22+34=' = '44+20='\n  '31+8=' = '90-81='\n  '62-27=' = '67-53='\n  '19+36=' = '69-7='\n  '99-10=' = '81-74='\n  '41+44=' = '11+65='\n  '78-63=' = '58-14='\n  '83-62=' = '21+50='\n  '39-32=' = '41-36='\n  '43+19=' = '74-0='\n  '3+21=' = '65-10='\n  '46-22=' = '26-24='\n  '29-16=' = '41+57='\n  '34-21=' = '90-41='\n  '86-62=' = '86-29='\n  '78+10=' = '83-7='\n  '63-31=' = '20+38='\n  '41+38=' = '47-29='\n  '94+1=' = '15+50='\n  '79+19=' = '85-44='\n  '14+15=' = '49+11='\n  '45-40=' = '6+50='\n  '86-14=' = '72-27='\n  '80-73=' = '22+39='\n  '44+35=' = '62-32='\n  '7+36=' = '99-88='\n  '91-61=' = '74-46='\n  '96-32=' = '65-55='\n  '59+16=' = '24+35='\n  '99-66=' = '27+16='\n  '63-51=' = '40+52='\n  '6+75=' = '49+8='\n  '19+61=' = '61+26='\n  '92-17=' = '43-3='\n  '24-20=' = '52-16='\n  '53-20=' = '99-14='\n  '39+40=' = '53-14='\n  '40+23=' = '24+59='\n  '0+52=' = '21-3='\n}\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $range = $cell.Range\n    # Cell text includes a trailing end-of-cell mark (cr + bell); trim it off.\n    $raw = $range.Text\n    $key = $raw.TrimEnd([char]13, [char]7)\n    if ($replacements.ContainsKey($key)) {\n      $range.Text = $replacements[$key]\n    }\n  }\n}\n"}
